# Update "南宁-漫展信息" workbook to the next scrape snapshot.
#
# Both the "展览" (sheet 1) and "全部类型" (sheet 4) tabs drop their oldest
# event (2024-07-27 良牙夏典, which has already happened / sold out) from the
# top of the list. Deleting that row shifts every following row up by one,
# which also drops the table's trailing row (the sheets shrink from
# A1:I8 -> A1:I7 and A1:I11 -> A1:I10 respectively) and keeps all the other
# event data intact. A handful of "want to go" counters (column F) also
# ticked up for events that are still upcoming in both sheets.

$wb = $excel.ActiveWorkbook

# want-to-go (column F) deltas for events that persist across the refresh,
# keyed by the bilibili show id in column H so they apply wherever the row
# ends up after the shift.
$wantToGoUpdates = @{
    'id=86994' = 358   # 南宁·火影忍者only            357 -> 358
    'id=85370' = 752   # 南宁·蔚蓝档案only            748 -> 752
    'id=88227' = 259   # 南宁·国乙only                257 -> 259
    'id=89145' = 745   # 南宁·熊喵M动漫嘉年华【免费】  735 -> 745
    'id=88276' = 1839  # 南宁·第二届北极光动漫展      1811 -> 1839
}

foreach ($sheetIndex in 1, 4) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    # Drop the obsolete first event row; everything below shifts up and the
    # sheet's dimension/used-range shrinks by one row automatically.
    $ws.Rows.Item(2).Delete()

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

    for ($r = 2; $r -le $lastRow; $r++) {
        # Column A is a 0-based running index (header row = 0), renumber it
        # now that row 2's original occupant is gone.
        $ws.Cells.Item($r, 1).Value = $r - 1

        # Refresh the "want to go" counter for whichever rows matched.
        # (.Value2 is used for the read-back here - .Value doesn't surface a
        # usable scalar for comparisons in this host.)
        $link = $ws.Cells.Item($r, 8).Value2
        foreach ($id in $wantToGoUpdates.Keys) {
            if ($link -like "*$id*") {
                $ws.Cells.Item($r, 6).Value = $wantToGoUpdates[$id]
            }
        }
    }
}
